$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add three new header cells in row 1 (columns D, E, F)
$ws.Range("D1").Value = "ORG_LINE_IDENOLD"
$ws.Range("E1").Value = "ORG_LINE_IDENNEW"
$ws.Range("F1").Value = "ORG_LINE_STATUS"

# Update the window selection to D8
$ws.Range("D8").Select()
